# Add "Handy links" section (with hyperlinks) to slide 5's content placeholder,
# matching the commit "Added handy links in powerpoint".

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# The placeholder shrinks text to fit ("Shrink text on overflow") once the new
# links push it past the box bounds.
$tf.AutoSize = 2

# Track the 1-indexed insertion cursor; it starts right after the existing text.
$pos = $tr.Text.Length + 1

function Ins([string]$text) {
    $tr.InsertAfter($text) | Out-Null
    $len = $text.Length
    $range = $tr.Characters($pos, $len)
    $pos = $pos + $len
    return $range
}

# New blank paragraph separating the existing content from the links section.
Ins("`r") | Out-Null

# --- "Handy links (Not necessarily needed for this course):" (bold) ---
Ins("`r") | Out-Null
$r = Ins("Handy links ("); $r.Font.Bold = $true
$r = Ins("Not"); $r.Font.Bold = $true
$r = Ins(" "); $r.Font.Bold = $true
$r = Ins("necessarily"); $r.Font.Bold = $true
$r = Ins(" "); $r.Font.Bold = $true
$r = Ins("needed"); $r.Font.Bold = $true
$r = Ins(" "); $r.Font.Bold = $true
$r = Ins("for"); $r.Font.Bold = $true
$r = Ins(" "); $r.Font.Bold = $true
$r = Ins("this"); $r.Font.Bold = $true
$r = Ins(" course):"); $r.Font.Bold = $true

# --- "Awesome Polymer: https://github.com/Granze/awesome-polymer" ---
Ins("`r") | Out-Null
Ins("Awesome Polymer: ") | Out-Null
$link = "https://github.com/Granze/awesome-polymer"
$r = Ins($link)
$r.Font.Underline = $true
$r.ActionSettings.Item(1).Hyperlink.Address = $link

# --- "Handling events in Polymer: https://alligator.io/polymer/handling-events/" ---
Ins("`r") | Out-Null
Ins("Handling events in Polymer: ") | Out-Null
$link = "https://alligator.io/polymer/handling-events/"
$r = Ins($link)
$r.Font.Underline = $true
$r.ActionSettings.Item(1).Hyperlink.Address = $link

# --- "Polymer 2 cheatsheet: https://meowni.ca/posts/polymer-2-cheatsheet/" ---
Ins("`r") | Out-Null
Ins("Polymer 2 ") | Out-Null
Ins("cheatsheet") | Out-Null
Ins(": ") | Out-Null
$link = "https://meowni.ca/posts/polymer-2-cheatsheet/"
$r = Ins($link)
$r.Font.Underline = $true
$r.ActionSettings.Item(1).Hyperlink.Address = $link

# --- "Web components library: https://www.webcomponents.org" ---
Ins("`r") | Out-Null
Ins("Web components library: ") | Out-Null
$link = "https://www.webcomponents.org"
$r = Ins($link)
$r.Font.Underline = $true
$r.ActionSettings.Item(1).Hyperlink.Address = $link
